# BP terminal gate pricing - roll the "Daily TGP" table forward by one day.
# For every state/terminal block the newest-date rows are pushed down into the
# previous-newest slot (their values are unchanged) and a brand new top block
# of prices is written in for the new effective date. Column layout is:
#   A = Effective Date, B = (blank spacer), C = Terminal, D = Diesel,
#   E = ULP, F = PULP, G = e10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-Row {
    param(
        [int]$Row,
        [double]$Date,
        $Diesel,
        $ULP,
        $PULP,
        $E10
    )

    $ws.Cells.Item($Row, 1).Value2 = $Date
    if ($null -ne $Diesel) { $ws.Cells.Item($Row, 4).Value2 = $Diesel }
    if ($null -ne $ULP)    { $ws.Cells.Item($Row, 5).Value2 = $ULP }
    if ($null -ne $PULP)   { $ws.Cells.Item($Row, 6).Value2 = $PULP }
    if ($null -ne $E10)    { $ws.Cells.Item($Row, 7).Value2 = $E10 }
}

# --- New South Wales (Sydney-Botany, Sydney-Silverwater, Newcastle) ---
Set-Row  8 45993 164.35 160.5   170.5   160.61000000000001
Set-Row  9 45993 164.35 160.5   170.5   160.61000000000001
Set-Row 10 45993 165.96 163.34  173.34  163.86
Set-Row 11 45990 165.46 160.34  170.34  160.5
Set-Row 12 45990 165.46 160.34  170.34  160.5
Set-Row 13 45990 166.84 162.85  172.85  163.34

# --- Northern Territory (Darwin) ---
Set-Row 17 45993 169.08 165.68 175.68
Set-Row 18 45990 170.02 165.22 175.22

# --- Queensland (Brisbane, Cairns, Gladstone, Mackay, Townsville) ---
Set-Row 22 45993 165.37 162.57  172.17  163.72999999999999
Set-Row 23 45993 170.75 167.36  177.36
Set-Row 24 45993 170.52 167.81  177.81
Set-Row 25 45993 171.35 167.23  177.23  167
Set-Row 26 45993 169.99 168.69  178.69
Set-Row 27 45990 166.27 162.09  171.69  163.38
Set-Row 28 45990 171.63 166.87  176.87
Set-Row 29 45990 171.41 167.3   177.3
Set-Row 30 45990 172.24 166.71  176.71  166.75
Set-Row 31 45990 170.89 168.16  178.16

# --- South Australia (Adelaide) ---
Set-Row 35 45993 164.32 160.33000000000001 169.33
Set-Row 36 45990 165.2  159.84             168.84

# --- Tasmania (Burnie, Hobart) ---
Set-Row 40 45993 170.1  165.96 175.96
Set-Row 41 45993 169.81 166.38 176.38
Set-Row 42 45990 170.94 165.4  175.4
Set-Row 43 45990 170.65 165.82 175.82

# --- Victoria (Geelong, Melbourne) ---
Set-Row 47 45993 165.87 161.37 171.37
Set-Row 48 45993 165.77 161.49 171.49
Set-Row 49 45990 168.41 161.24 171.24
Set-Row 50 45990 168.33 161.37 171.37

# --- Western Australia (Broome, Esperance, Geraldton, Kalgoorlie, Perth, Port Hedland) ---
Set-Row 54 45993 179.93 177.34 187.34
Set-Row 55 45993 167.93 173.28 183.28
Set-Row 56 45993 170.26
Set-Row 57 45993 169.48 167.55
Set-Row 58 45993 165.39 163.6  173.6
Set-Row 59 45993 172    174.98
Set-Row 60 45990 180.86 176.62 186.62
Set-Row 61 45990 168.81 172.85 182.85
Set-Row 62 45990 171.13
Set-Row 63 45990 170.43 167.11
Set-Row 64 45990 166.34 163.16999999999999 173.17
Set-Row 65 45990 172.95 174.36
